$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift name/age data from columns B/C into columns A/B, drop the old
# numeric index column, and append the two new rows.
$ws.Range("A1").Value = "vinoth"
$ws.Range("B1").Value = 28
$ws.Range("A2").Value = "vimal"
$ws.Range("B2").Value = 26
$ws.Range("A3").Value = "sridhar"
$ws.Range("B3").Value = 59
$ws.Range("A4").Value = "vijaya"
$ws.Range("B4").Value = 55

# Clear out the old column C (ages used to live there) now that it is unused.
$ws.Range("C1:C2").Clear()

$ws.Range("B3").Select()
